$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value -replace "_old$", "_FV2410")
}
foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value -replace "_new$", "_FV2504")
}

# --- 2. Turn the used range into an Excel Table (ListObject) named "Table1" ---
$rng = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1, $null)
$lo.Name = "Table1"

# --- 3. Freeze the header row (split/freeze pane after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
